# Apply "alter notebook new chart" edit to estrategias worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# A1 previously held the "Estratégia" header label; it is removed entirely.
$ws.Range("A1").Clear()

# B1:K1 previously held text headers (Related, Visited, ...); they become
# plain numeric index values 0-9, keeping their existing style.
$headerVals = @(0,1,2,3,4,5,6,7,8,9)
for ($i = 0; $i -lt $headerVals.Length; $i++) {
    $col = 2 + $i   # B=2 .. K=11
    $ws.Cells.Item(1, $col).Value = $headerVals[$i]
}

# --- Data rows (rows 2-8), column A becomes numeric 0-6, column B keeps the
#     "E1".."E7" labels (now referencing the compacted shared string table),
#     and columns F, G, H, I, K get updated to full-precision computed values.
$rowData = @(
    @{ Row = 2; Idx = 0; Label = "E1"; F = 0.04426559356136821; G = 0.4313725490196079;  H = 0.08029197080291971; I = 0.03018108651911469; J = 0.5;                K = 0.05692599620493359 },
    @{ Row = 3; Idx = 1; Label = "E2"; F = 0.02722904431393487; G = 1;                    H = 0.05301455301455301; I = 0.01601708489054992; J = 1;                  K = 0.03152916447714135 },
    @{ Row = 4; Idx = 2; Label = "E3"; F = 0.0339943342776204;  G = 0.7058823529411765;   H = 0.06486486486486487; I = 0.02077431539187913; J = 0.7333333333333333; K = 0.04040404040404041 },
    @{ Row = 5; Idx = 3; Label = "E4"; F = 0.03747870528109029; G = 0.8627450980392157;   H = 0.07183673469387755; I = 0.02129471890971039; J = 0.8333333333333334; K = 0.04152823920265781 },
    @{ Row = 6; Idx = 4; Label = "E5"; F = 0.06506849315068493; G = 0.3725490196078431;   H = 0.1107871720116618;  I = 0.04794520547945205; J = 0.4666666666666667; K = 0.08695652173913043 },
    @{ Row = 7; Idx = 5; Label = "E6"; F = 0.06194690265486726; G = 0.6862745098039216;   H = 0.1136363636363636;  I = 0.03539823008849557; J = 0.6666666666666666; K = 0.06722689075630252 },
    @{ Row = 8; Idx = 6; Label = "E7"; F = 0.05811138014527845; G = 0.4705882352941176;   H = 0.103448275862069;   I = 0.03631961259079903; J = 0.5;                K = 0.06772009029345372 }
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Idx     # column A -> numeric index
    $ws.Cells.Item($r, 2).Value = $entry.Label   # column B -> "E1".."E7"
    $ws.Cells.Item($r, 6).Value = $entry.F        # column F -> Precision
    $ws.Cells.Item($r, 7).Value = $entry.G        # column G -> Recall
    $ws.Cells.Item($r, 8).Value = $entry.H        # column H -> F-Measure
    $ws.Cells.Item($r, 9).Value = $entry.I        # column I -> Final Precision
    $ws.Cells.Item($r, 10).Value = $entry.J       # column J -> Final Recall
    $ws.Cells.Item($r, 11).Value = $entry.K       # column K -> Final F-Measure
}
